$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume snapshot (column D = Price, column E =
# Volume(1h)) plus two rows whose coin entries were reordered/replaced.
# Column D values are stored as plain text in this sheet (e.g. "26.498.68",
# "0.07800", "1.000"), so a leading apostrophe is used to force Excel to keep
# them as text instead of re-interpreting them as numbers/dates and silently
# dropping significant trailing/leading zeros.
$ws.Range("D2").Value = "'26.498.68"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "'1.671.32"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'220.02"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").Value = "'0.5273"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.2676"
$ws.Range("E8").Value = "  +2.35%  "
$ws.Range("D9").Value = "'0.06374"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "'21.74"
$ws.Range("E10").Value = "  +3.94%  "
$ws.Range("D11").Value = "'0.07800"
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").Value = "'1.673.50"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "'4.487"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "'0.5563"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'0.0₅8285"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "'65.59"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "'26.510.48"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'4.759"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").Value = "'193.08"
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("D21").Value = "'10.33"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'6.304"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'0.1268"
$ws.Range("E24").Value = "  +3.80%  "
$ws.Range("D25").Value = "'138.25"
$ws.Range("E25").Value = "  -5.36%  "
$ws.Range("D26").Value = "'7.392"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'16.32"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("E28").Value = "  +2.32%  "
$ws.Range("D29").Value = "'0.06237"
$ws.Range("E29").Value = "  +4.56%  "
$ws.Range("D30").Value = "'1.289"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("E31").Value = "  +6.01%  "
$ws.Range("D32").Value = "'3.416"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "'1.689"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "'1.007"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").Value = "'0.6159"
$ws.Range("E35").Value = "  +9.05%  "
$ws.Range("D36").Value = "'2.422"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("D37").Value = "'2.786"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").Value = "'0.01616"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").Value = "'6.040"
$ws.Range("E39").Value = "  +3.39%  "
$ws.Range("D40").Value = "'1.093.40"
$ws.Range("E40").Value = "  +6.38%  "
$ws.Range("D41").Value = "'0.8588"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").Value = "'1.816.94"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "'58.63"
$ws.Range("E45").Value = "  +5.01%  "
$ws.Range("D46").Value = "'0.0₈106"
$ws.Range("E46").Value = "  -4.87%  "
$ws.Range("D47").Value = "'8.198"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.000"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.516"
$ws.Range("E49").Value = "  +9.71%  "
$ws.Range("D50").Value = "'0.05197"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.4234"
$ws.Range("E51").Value = "  +0.42%  "
